$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 17, shifting existing rows 17-82 down to 18-83
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with a new weekly record
$ws.Cells.Item(17, 1).Value = 2
$ws.Cells.Item(17, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44910
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100112026
$ws.Cells.Item(17, 7).Value = "Haba"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 500
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 8000
$ws.Cells.Item(17, 13).Value = 7500
$ws.Cells.Item(17, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 16).Value = 300
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Match date style/format used by the rest of column D
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
